$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A140").Value = 139
$ws.Range("B140").Value = 1
$ws.Range("C140").Value = "2024-06-17 21:12:19"
$ws.Range("D140").Value = 200
$ws.Range("E140").Value = 17

$ws.Range("A141").Value = 140
$ws.Range("B141").Value = 2
$ws.Range("C141").Value = "2024-06-17 21:12:20"
$ws.Range("D141").Value = 200
$ws.Range("E141").Value = 1
